$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update the existing product's data
$ws.Range("A2").Value = "Nettoyant multi-surfaces à la menthe et à l'eucalyptus"
$ws.Range("B2").Value = 7.8
$ws.Range("C2").Value = "Spray de 750 ml
Poids net 0.750000
Référence 1333803
CONSEILS D’UTILISATION
Pour nettoyer et faire briller les sanitaires ( robinetterie, lavabos, carrelage) vaporiser sur la surface, laisser agir, frotter si besoin, rincer puis essuyer. Pour respecter l'environnement, suivre les dosages préconisés. Pour plus d'informations, www.biovie.com.
COMPOSITION
AQUA,CITRIC ACID,LACTIC ACID,CAPRYL GLUCOSIDE,AMMONIUM LAURYL SULFATE,ALCOHOL,ALCOHOL C12-18 ETHOXYLATED,SODIUM CITRATE,ORGANIC MENTHA VIRIDIS LEAF OIL,PARFUM,ISOPROPYL ALCOHOL,LIMONENE,LINALOOL,ORGANIC EUCALYPTUS GLOBULUS OIL"
$ws.Range("D2").Value = 30
$ws.Range("E2").Value = "produit d'entretient"

# Row 3: add a second product
$ws.Range("A3").Value = "Savon noir liquide à l'huile de lin"
$ws.Range("B3").Value = 4.5
$ws.Range("C3").Value = "Bouteille d'1L
Poids net 1.026000
Référence 1333080
CONSEILS D’UTILISATION
- Cuisine : Mettre un peu de savon noir sur une éponge. Laver et rincer à l'eau chaude.
- Sol carrelé : Diluer 2-3 bouchons dans un seau d'eau chaude.
- Linge :  Étaler un peu de savon noir sur la tache. Frictionner délicatement et placer en machine. Toujours faire un essai préalable. Peut être utilisé comme lessive pour linge délicat. 4 à 6 bouchons par machine.
Précautions : provoque une sévère irritation des yeux. Tenir hors de portée des enfants. En cas de consultation d’un médecin, garder à disposition le récipient ou l’étiquette. Porter un équipement de protection des yeux, un équipement de protection du visage.
En cas de contact avec les yeux : rincer avec précaution à l’eau pendant plusieurs minutes. Enlever les lentilles de contact si la victime en porte et si elles peuvent être facilement enlevées. Continuer à rincer.
Si l’irritation oculaire persiste: consulter un médecin. Aérer les pièces au moins dix minutes pendant et après le nettoyage, en été comme en hiver.
COMPOSITION
AQUA FATTY ACIDS, C16-18 AND C18-UNSATD, POTASSIUM SALTS,POTASSIUM CARBONATE,GLYCERIN,SODIUM CITRATE,TRISODIUM SALT OF METHYL GLYCINEDIACETIC ACID,LINUM USITATISSIMUM OIL"
$ws.Range("D3").Value = 45
$ws.Range("E3").Value = "produit d'entretient"

# Column C needs to be wide enough for the long description text
$ws.Columns.Item(3).ColumnWidth = 48.166667

# Row 3 wraps to its natural auto height (no explicit custom height)
$ws.Rows.Item(3).AutoFit()

# Row 2 keeps an explicit custom height
$ws.Rows.Item(2).RowHeight = 26.5

$null = $ws.Range("E8").Select()
